$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("I1:J1")
$rng.FormulaArray = '=CELL("width")'

$ws.Range("I2").Select()
